# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Goblin_Profits workbook. For each worksheet, specific cells in the
# currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# columns are updated to reflect newly scraped market board values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3000
$ws.Range("I74").Value = 3000
$ws.Range("K74").Value = 3000
$ws.Range("M74").Value = -2064
$ws.Range("H77").Value = 3000
$ws.Range("I77").Value = 3000
$ws.Range("K77").Value = 15000
$ws.Range("M77").Value = -10320
$ws.Range("H103").Value = 1650
$ws.Range("J103").Value = 1733.3334
$ws.Range("L103").Value = 5200.0002
$ws.Range("N103").Value = -6372.0002
$ws.Range("H116").Value = 6923.241
$ws.Range("I116").Value = 7165.8335
$ws.Range("K116").Value = 7165.8335
$ws.Range("M116").Value = -3723.8335
$ws.Range("H133").Value = 93247.625
$ws.Range("J133").Value = 93247.625
$ws.Range("L133").Value = 93247.625
$ws.Range("N133").Value = -103367.625
$ws.Range("H134").Value = 291427.84
$ws.Range("J134").Value = 291427.84
$ws.Range("L134").Value = 291427.84
$ws.Range("N134").Value = -301567.84
$ws.Range("H137").Value = 1835.2759
$ws.Range("J137").Value = 2143.4
$ws.Range("L137").Value = 6430.200000000001
$ws.Range("N137").Value = -11530.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1728.1818
$ws.Range("I45").Value = 1412
$ws.Range("K45").Value = 1412
$ws.Range("M45").Value = -1035
$ws.Range("H63").Value = 4396.1333
$ws.Range("I63").Value = 1326.8889
$ws.Range("J63").Value = 9000
$ws.Range("K63").Value = 1326.8889
$ws.Range("L63").Value = 9000
$ws.Range("M63").Value = -640.8888999999999
$ws.Range("N63").Value = -10372
$ws.Range("H66").Value = 4396.1333
$ws.Range("I66").Value = 1326.8889
$ws.Range("J66").Value = 9000
$ws.Range("K66").Value = 6634.4445
$ws.Range("L66").Value = 45000
$ws.Range("M66").Value = -3202.4445
$ws.Range("N66").Value = -51864

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 24399.8
$ws.Range("I26").Value = 24399.8
$ws.Range("K26").Value = 24399.8
$ws.Range("M26").Value = -24107.8
$ws.Range("H86").Value = 5304.8076
$ws.Range("I86").Value = 5788.0435
$ws.Range("K86").Value = 5788.0435
$ws.Range("M86").Value = -4665.0435
$ws.Range("H89").Value = 5304.8076
$ws.Range("I89").Value = 5788.0435
$ws.Range("K89").Value = 28940.2175
$ws.Range("M89").Value = -23324.2175
$ws.Range("H105").Value = 4867.2144
$ws.Range("I105").Value = 6073.074
$ws.Range("J105").Value = 2696.6667
$ws.Range("K105").Value = 6073.074
$ws.Range("L105").Value = 2696.6667
$ws.Range("M105").Value = -4326.074
$ws.Range("N105").Value = -6190.6667
$ws.Range("H135").Value = 261249.75
$ws.Range("I135").Value = 94999
$ws.Range("J135").Value = 316666.66
$ws.Range("K135").Value = 94999
$ws.Range("L135").Value = 316666.66
$ws.Range("M135").Value = -89929
$ws.Range("N135").Value = -326806.66

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 25788.3
$ws.Range("I50").Value = 19966.666
$ws.Range("K50").Value = 19966.666
$ws.Range("M50").Value = -19341.666
$ws.Range("H58").Value = 2515.7827
$ws.Range("I58").Value = 2521.0454
$ws.Range("K58").Value = 2521.0454
$ws.Range("M58").Value = -2318.0454
$ws.Range("H93").Value = 10116.454
$ws.Range("I93").Value = 8128.1
$ws.Range("K93").Value = 8128.1
$ws.Range("M93").Value = -6256.1
$ws.Range("H94").Value = 2370
$ws.Range("I94").Value = 1200
$ws.Range("J94").Value = 2537.1428
$ws.Range("K94").Value = 1200
$ws.Range("L94").Value = 2537.1428
$ws.Range("M94").Value = -749
$ws.Range("N94").Value = -3439.1428
$ws.Range("H99").Value = 2313.875
$ws.Range("I99").Value = 2306.1667
$ws.Range("K99").Value = 2306.1667
$ws.Range("M99").Value = -808.1667000000002
$ws.Range("H107").Value = 838.5
$ws.Range("I107").Value = 922.1111
$ws.Range("K107").Value = 922.1111
$ws.Range("M107").Value = 997.8889
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("L124").Value = 0
$ws.Range("H126").Value = 2313.875
$ws.Range("I126").Value = 2306.1667
$ws.Range("K126").Value = 6918.500100000001
$ws.Range("M126").Value = -4448.500100000001
$ws.Range("H136").Value = 2515.7827
$ws.Range("I136").Value = 2521.0454
$ws.Range("K136").Value = 7563.1362
$ws.Range("M136").Value = -5013.1362
$ws.Range("H141").Value = 268068.3
$ws.Range("J141").Value = 268068.3
$ws.Range("L141").Value = 268068.3
$ws.Range("N141").Value = -278428.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3298.5
$ws.Range("J68").Value = 3148.9167
$ws.Range("L68").Value = 9446.750100000001
$ws.Range("N68").Value = -11068.7501
$ws.Range("H71").Value = 3298.5
$ws.Range("J71").Value = 3148.9167
$ws.Range("L71").Value = 28340.2503
$ws.Range("N71").Value = -36452.2503
$ws.Range("H75").Value = 1000
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 1000
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H107").Value = 1333.1052
$ws.Range("I107").Value = 1903.3334
$ws.Range("J107").Value = 1069.9231
$ws.Range("K107").Value = 5710.0002
$ws.Range("L107").Value = 3209.7693
$ws.Range("M107").Value = -3790.0002
$ws.Range("N107").Value = -7049.7693
$ws.Range("H121").Value = 5216.8887
$ws.Range("I121").Value = 715
$ws.Range("J121").Value = 6503.143
$ws.Range("K121").Value = 2145
$ws.Range("L121").Value = 19509.429
$ws.Range("M121").Value = -835
$ws.Range("N121").Value = -22129.429
$ws.Range("H131").Value = 3929010
$ws.Range("J131").Value = 4452652
$ws.Range("L131").Value = 13357956
$ws.Range("N131").Value = -13368036

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 250038780
$ws.Range("J34").Value = 250038780
$ws.Range("L34").Value = 250038780
$ws.Range("N34").Value = -250039316
$ws.Range("H45").Value = 76666
$ws.Range("J45").Value = 76666
$ws.Range("L45").Value = 76666
$ws.Range("N45").Value = -77784
$ws.Range("H46").Value = 4500
$ws.Range("I46").Value = 4500
$ws.Range("K46").Value = 4500
$ws.Range("M46").Value = -4344
$ws.Range("H76").Value = 250038780
$ws.Range("J76").Value = 250038780
$ws.Range("L76").Value = 250038780
$ws.Range("N76").Value = -250039410
$ws.Range("H79").Value = 250038780
$ws.Range("J79").Value = 250038780
$ws.Range("L79").Value = 250038780
$ws.Range("N79").Value = -250040964
$ws.Range("H92").Value = 21750
$ws.Range("J92").Value = 21750
$ws.Range("L92").Value = 21750
$ws.Range("N92").Value = -25494
$ws.Range("H104").Value = 43447.332
$ws.Range("J104").Value = 43447.332
$ws.Range("L104").Value = 43447.332
$ws.Range("N104").Value = -50435.332
$ws.Range("H110").Value = 123750
$ws.Range("J110").Value = 123750
$ws.Range("L110").Value = 123750
$ws.Range("N110").Value = -131930
$ws.Range("H111").Value = 5000
$ws.Range("J111").Value = 5000
$ws.Range("L111").Value = 5000
$ws.Range("N111").Value = -11134
$ws.Range("H113").Value = 5629.237
$ws.Range("I113").Value = 2677.6
$ws.Range("J113").Value = 7554.2173
$ws.Range("K113").Value = 2677.6
$ws.Range("L113").Value = 7554.2173
$ws.Range("M113").Value = -507.5999999999999
$ws.Range("N113").Value = -11894.2173
$ws.Range("H114").Value = 35000
$ws.Range("J114").Value = 35000
$ws.Range("L114").Value = 35000
$ws.Range("N114").Value = -43678
$ws.Range("H120").Value = 18666.334
$ws.Range("J120").Value = 18666.334
$ws.Range("L120").Value = 18666.334
$ws.Range("N120").Value = -28342.334
$ws.Range("H126").Value = 3595
$ws.Range("I126").Value = 2743.75
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 8231.25
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -5761.25
$ws.Range("N126").Value = -25940
$ws.Range("H132").Value = 2558.9
$ws.Range("I132").Value = 2165.4546
$ws.Range("K132").Value = 6496.3638
$ws.Range("M132").Value = -3966.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2387.48
$ws.Range("I82").Value = 812.5
$ws.Range("J82").Value = 3841.3076
$ws.Range("K82").Value = 812.5
$ws.Range("L82").Value = 3841.3076
$ws.Range("M82").Value = -451.5
$ws.Range("N82").Value = -4563.3076
$ws.Range("H85").Value = 2387.48
$ws.Range("I85").Value = 812.5
$ws.Range("J85").Value = 3841.3076
$ws.Range("K85").Value = 812.5
$ws.Range("L85").Value = 3841.3076
$ws.Range("M85").Value = 435.5
$ws.Range("N85").Value = -6337.3076
$ws.Range("H133").Value = 98999
$ws.Range("J133").Value = 98999
$ws.Range("L133").Value = 98999
$ws.Range("N133").Value = -104059

$wb.Save()
Write-Host "Applied scheduled market data refresh to all sheets."